$wb = $excel.ActiveWorkbook

# Worksheet handles
$wsProperty        = $wb.Worksheets.Item("Property")
$wsBuildingList     = $wb.Worksheets.Item("Record_BuildingList")
$wsBuildingProduce  = $wb.Worksheets.Item("Record_BuildingProduce")

# --- Record_BuildingList: bump the SLG building config value (C2: 6 -> 8)
$wsBuildingList.Range("C2").Value = 8

# --- Property: add a new "LoadPropertyFinish" row (row 13)
$wsProperty.Range("A13").Value = "LoadPropertyFinish"
$wsProperty.Range("B13").Value = "int"
$wsProperty.Range("B13").NumberFormat = "@"
$wsProperty.Range("C13").Value = $true
$wsProperty.Range("D13").Value = $true
$wsProperty.Range("E13").Value = $true
$wsProperty.Range("G13").Value = 0
$wsProperty.Range("H13").Value = 0
$wsProperty.Range("I13").Value = "Friend"
$wsProperty.Range("I13").NumberFormat = "@"

# Extend the TRUE/FALSE list validation on column F down to cover the new row
$fValidationRange = $wsProperty.Range("F2:F1048576")
$fValidationRange.Validation.Delete()
$fValidationRange.Validation.Add(3, 1, 1, '"TRUE,FALSE"')

# --- Selection / active-sheet bookkeeping (order matters: last Select() wins the active tab)
$wsBuildingProduce.Range("F28").Select() | Out-Null
$wsProperty.Range("C26").Select() | Out-Null
$wsBuildingList.Range("G10").Select() | Out-Null
